# Add two new columns, "I0" (col I) and "IF" (col J), mirroring the
# existing header formatting (bold font, thin border, centered alignment)
# taken from the adjacent "IP" header cell (H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy H1's formatting onto the new header cells so they reuse the same
# style (rather than Excel inventing a new cellXfs entry).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data row values
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7
